$d = $word.ActiveDocument

# Locate the "ZEP" Heading 2 paragraph, then delete the very next paragraph
# (a leftover placeholder that contains only an italic "Zephaniah" run).
$paras = $d.Paragraphs
$count = $paras.Count
$targetIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $text = $p.Range.Text
    if ($text -eq "ZEP`r" -or $text -eq "ZEP") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $next = $paras.Item($targetIndex + 1)
    $nextText = $next.Range.Text
    if ($nextText -eq "Zephaniah`r" -or $nextText -eq "Zephaniah") {
        $next.Range.Delete()
    }
}
